$d = $word.ActiveDocument

$replacements = @(
    @("2025-09-15 Monday", "2025-09-16 Tuesday"),
    @("70×31=", "40×21="),
    @("48×49=", "12×85="),
    @("31×52=", "11×53="),
    @("91×83=", "44×88="),
    @("32×40=", "74×11="),
    @("50×89=", "50×45="),
    @("62×57=", "71×57="),
    @("25×35=", "34×78="),
    @("53×18=", "33×15="),
    @("21×31=", "87×76="),
    @("49×23=", "65×37="),
    @("89×85=", "32×92="),
    @("61×80=", "55×97="),
    @("13×36=", "25×31="),
    @("60×11=", "15×59="),
    @("95×30=", "14×96="),
    @("11×89=", "64×41="),
    @("30×34=", "60×37="),
    @("25×83=", "28×86="),
    @("22×82=", "96×22="),
    @("90×57=", "95×32="),
    @("19×37=", "68×97="),
    @("38×23=", "13×37="),
    @("12×34=", "57×54="),
    @("57×75=", "86×93=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
